# Applies the editorial-tag-markup revisions described by the diff:
#   - splits plain runs so certain words are wrapped in <ms>...</ms>,
#     <bp>...</bp>, <tmp>...</tmp> markers (rendered in blue Courier New)
#   - wraps an existing <m>...</m> markup pair in <tl>...</tl>

$d = $word.ActiveDocument

function Set-TagFont($range) {
    $range.Font.Name  = "Courier New"
    $range.Font.Size  = 9
    $range.Font.Color = 16711680   # BGR encoding of RGB 0000ff (blue)
}

# ---------------------------------------------------------------------
# Change 1: " Et prandre trois onces " ->
#   " Et prandre trois " <ms> "onces" </ms> " "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" Et prandre trois onces ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origStart = $rng.Start
$newText = " Et prandre trois <ms>onces</ms> "
$rng.Text = $newText

$openStart = $origStart + $newText.IndexOf("<ms>")
$openEnd = $openStart + 4
Set-TagFont ($d.Range($openStart, $openEnd))

$closeStart = $origStart + $newText.IndexOf("</ms>")
$closeEnd = $closeStart + 5
Set-TagFont ($d.Range($closeStart, $closeEnd))

# ---------------------------------------------------------------------
# Change 2: " choppine de" -> " " <ms> "choppine" </ms> " de"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" choppine de", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origStart = $rng.Start
$newText = " <ms>choppine</ms> de"
$rng.Text = $newText

$openStart = $origStart + $newText.IndexOf("<ms>")
$openEnd = $openStart + 4
Set-TagFont ($d.Range($openStart, $openEnd))

$closeStart = $origStart + $newText.IndexOf("</ms>")
$closeEnd = $closeStart + 5
Set-TagFont ($d.Range($closeStart, $closeEnd))

# ---------------------------------------------------------------------
# Change 3: " Et faire bouillir les matieres lespace de demye" ->
#   " Et faire bouillir les matieres l" <ms> "espace de " <tmp> "demye"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" Et faire bouillir les matieres lespace de demye", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origStart = $rng.Start
$newText = " Et faire bouillir les matieres l<ms>espace de <tmp>demye"
$rng.Text = $newText

$msStart = $origStart + $newText.IndexOf("<ms>")
$msEnd = $msStart + 4
Set-TagFont ($d.Range($msStart, $msEnd))

$tmpStart = $origStart + $newText.IndexOf("<tmp>")
$tmpEnd = $tmpStart + 5
Set-TagFont ($d.Range($tmpStart, $tmpEnd))

# ---------------------------------------------------------------------
# Change 4: "heure &" -> "heure" </tmp></ms> " &"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("heure &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origStart = $rng.Start
$newText = "heure</tmp></ms> &"
$rng.Text = $newText

$closeStart = $origStart + $newText.IndexOf("</tmp></ms>")
$closeEnd = $closeStart + 11
Set-TagFont ($d.Range($closeStart, $closeEnd))

# ---------------------------------------------------------------------
# Change 5: <m>pappier</m> -> <tl><m>pappier</m></tl>
# (edits the text of the two existing marker runs in place, preserving
#  their original character formatting)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("<m>pappier", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mOpenStart = $rng.Start
$d.Range($mOpenStart, $mOpenStart + 3).Text = "<tl><m>"

$rng = $d.Content
$rng.Find.Execute("pappier</m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mCloseStart = $rng.Start + 7
$mCloseEnd = $rng.End
$d.Range($mCloseStart, $mCloseEnd).Text = "</m></tl>"

# ---------------------------------------------------------------------
# Change 6: "poinct la main en vos matieres" ->
#   "poinct la " <bp> "main" </bp> " en vos matieres"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("poinct la main en vos matieres", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origStart = $rng.Start
$newText = "poinct la <bp>main</bp> en vos matieres"
$rng.Text = $newText

$openStart = $origStart + $newText.IndexOf("<bp>")
$openEnd = $openStart + 4
Set-TagFont ($d.Range($openStart, $openEnd))

$closeStart = $origStart + $newText.IndexOf("</bp>")
$closeEnd = $closeStart + 5
Set-TagFont ($d.Range($closeStart, $closeEnd))

Write-Host "done"
